# Update handback status report timestamps for the latest handback run.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 3 corresponds to ee06cae0...zh-cn.xlf
$wsZhCn.Range("E3").Value = "2016-03-11 16:31:55"
$wsZhCn.Range("H3").Value = "2016-03-11 16:32:13"

# de-de sheet: row 3 corresponds to ee06cae0...de-de.xlf
$wsDeDe.Range("E3").Value = "2016-03-11 16:31:58"
$wsDeDe.Range("H3").Value = "2016-03-11 16:32:18"
